$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellValue {
    param($Ws, $Row, $Col, $Val)
    if ($Val -match '^-?\d+(\.\d+)?$') {
        $Ws.Cells.Item($Row, $Col).Value = "'" + $Val
    } else {
        $Ws.Cells.Item($Row, $Col).Value = $Val
    }
}

Set-CellValue $ws 2 4 "26.380.61"
Set-CellValue $ws 2 5 "  -8.09%  "
Set-CellValue $ws 3 4 "1.678.40"
Set-CellValue $ws 3 5 "  -6.88%  "
Set-CellValue $ws 4 5 "  +0.45%  "
Set-CellValue $ws 5 4 "216.48"
Set-CellValue $ws 5 5 "  -6.57%  "
Set-CellValue $ws 6 4 "1.008"
Set-CellValue $ws 6 5 "  +0.36%  "
Set-CellValue $ws 7 4 "0.4946"
Set-CellValue $ws 7 5 "  -16.82%  "
Set-CellValue $ws 8 4 "0.2600"
Set-CellValue $ws 8 5 "  -6.36%  "
Set-CellValue $ws 9 4 "21.61"
Set-CellValue $ws 9 5 "  -7.40%  "
Set-CellValue $ws 10 4 "0.06120"
Set-CellValue $ws 10 5 "  -10.37%  "
Set-CellValue $ws 11 4 "0.07285"
Set-CellValue $ws 11 5 "  -3.40%  "
Set-CellValue $ws 12 4 "1.712.50"
Set-CellValue $ws 12 5 "  -4.82%  "
Set-CellValue $ws 13 4 "4.401"
Set-CellValue $ws 13 5 "  -8.09%  "
Set-CellValue $ws 14 4 "0.5708"
Set-CellValue $ws 15 4 "1.908.50"
Set-CellValue $ws 15 5 "  -6.82%  "
Set-CellValue $ws 16 4 "0.000008170"
Set-CellValue $ws 16 5 "  -11.99%  "
Set-CellValue $ws 17 4 "64.06"
Set-CellValue $ws 17 5 "  -14.98%  "
Set-CellValue $ws 18 4 "26.429.41"
Set-CellValue $ws 18 5 "  -7.89%  "
Set-CellValue $ws 19 4 "4.973"
Set-CellValue $ws 19 5 "  -9.09%  "
Set-CellValue $ws 20 5 "  +0.29%  "
Set-CellValue $ws 21 4 "10.68"
Set-CellValue $ws 21 5 "  -6.77%  "
Set-CellValue $ws 22 4 "181.90"
Set-CellValue $ws 22 5 "  -13.52%  "
Set-CellValue $ws 23 2 "BinanceUSD"
Set-CellValue $ws 23 3 "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-CellValue $ws 23 4 "1.008"
Set-CellValue $ws 23 5 "  +0.40%  "
Set-CellValue $ws 24 2 "Chainlink"
Set-CellValue $ws 24 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-CellValue $ws 24 4 "6.139"
Set-CellValue $ws 24 5 "  -10.46%  "
Set-CellValue $ws 25 4 "144.03"
Set-CellValue $ws 25 5 "  -6.65%  "
Set-CellValue $ws 26 4 "7.493"
Set-CellValue $ws 26 5 "  -4.52%  "
Set-CellValue $ws 27 4 "0.1124"
Set-CellValue $ws 27 5 "  -11.85%  "
Set-CellValue $ws 28 4 "15.43"
Set-CellValue $ws 28 5 "  -5.91%  "
Set-CellValue $ws 29 4 "1.309"
Set-CellValue $ws 29 5 "  -8.50%  "
Set-CellValue $ws 30 4 "0.05625"
Set-CellValue $ws 30 5 "  -8.87%  "
Set-CellValue $ws 31 4 "1.318"
Set-CellValue $ws 31 5 "  -7.21%  "
Set-CellValue $ws 32 4 "3.462"
Set-CellValue $ws 32 5 "  -8.47%  "
Set-CellValue $ws 33 4 "3.443"
Set-CellValue $ws 33 5 "  -8.15%  "
Set-CellValue $ws 34 4 "1.624"
Set-CellValue $ws 34 5 "  -5.45%  "
Set-CellValue $ws 35 4 "0.9999"
Set-CellValue $ws 35 5 "  -5.85%  "
Set-CellValue $ws 36 4 "2.368"
Set-CellValue $ws 36 5 "  -5.09%  "
Set-CellValue $ws 37 4 "0.5851"
Set-CellValue $ws 37 5 "  -8.63%  "
Set-CellValue $ws 38 4 "2.626"
Set-CellValue $ws 38 5 "  -3.26%  "
Set-CellValue $ws 39 4 "0.01576"
Set-CellValue $ws 39 5 "  -7.89%  "
Set-CellValue $ws 40 4 "1.068.70"
Set-CellValue $ws 40 5 "  -5.66%  "
Set-CellValue $ws 41 4 "5.895"
Set-CellValue $ws 41 5 "  -8.33%  "
Set-CellValue $ws 42 4 "0.8494"
Set-CellValue $ws 42 5 "  -2.07%  "
Set-CellValue $ws 43 4 "1.002"
Set-CellValue $ws 43 5 "  -0.17%  "
Set-CellValue $ws 44 4 "98.13"
Set-CellValue $ws 44 5 "  -2.45%  "
Set-CellValue $ws 45 4 "1.837.31"
Set-CellValue $ws 45 5 "  -6.32%  "
Set-CellValue $ws 46 4 "56.07"
Set-CellValue $ws 46 5 "  -7.46%  "
Set-CellValue $ws 47 4 "0.00000000105"
Set-CellValue $ws 47 5 "  -6.48%  "
Set-CellValue $ws 48 4 "1.000"
Set-CellValue $ws 48 5 "  -0.28%  "
Set-CellValue $ws 49 4 "8.022"
Set-CellValue $ws 49 5 "  -3.54%  "
Set-CellValue $ws 50 4 "0.4327"
Set-CellValue $ws 50 5 "  -3.69%  "
Set-CellValue $ws 51 4 "0.05180"
Set-CellValue $ws 51 5 "  -5.33%  "
